$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 347.81818
$ws.Cells.Item(15, 9).Value = 347.81818
$ws.Cells.Item(15, 11).Value = 1043.45454
$ws.Cells.Item(15, 13).Value = -874.45454
$ws.Cells.Item(40, 8).Value = 1699.5
$ws.Cells.Item(40, 9).Value = 1400.5
$ws.Cells.Item(40, 11).Value = 1400.5
$ws.Cells.Item(40, 13).Value = -1225.5
$ws.Cells.Item(54, 8).Value = 2499.5
$ws.Cells.Item(54, 10).Value = 0
$ws.Cells.Item(54, 12).Value = 0
$ws.Cells.Item(54, 14).ClearContents()
$ws.Cells.Item(96, 8).Value = 664.3333
$ws.Cells.Item(96, 9).Value = 598.2
$ws.Cells.Item(96, 11).Value = 1794.6
$ws.Cells.Item(96, 13).Value = -421.6000000000001
$ws.Cells.Item(100, 8).Value = 2359.5
$ws.Cells.Item(100, 9).Value = 2211.875
$ws.Cells.Item(100, 11).Value = 2211.875
$ws.Cells.Item(100, 13).Value = -1670.875
$ws.Cells.Item(111, 8).Value = 389.25
$ws.Cells.Item(111, 9).Value = 389.25
$ws.Cells.Item(111, 11).Value = 1167.75
$ws.Cells.Item(111, 13).Value = 1899.25
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 863.7273
$ws.Cells.Item(2, 9).Value = 800.1
$ws.Cells.Item(2, 11).Value = 800.1
$ws.Cells.Item(2, 13).Value = -687.1
$ws.Cells.Item(32, 8).Value = 1806.3549
$ws.Cells.Item(32, 9).Value = 1660.9323
$ws.Cells.Item(32, 11).Value = 1660.9323
$ws.Cells.Item(32, 13).Value = -1373.9323
$ws.Cells.Item(116, 8).Value = 863.7273
$ws.Cells.Item(116, 9).Value = 800.1
$ws.Cells.Item(116, 11).Value = 800.1
$ws.Cells.Item(116, 13).Value = 1493.9
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 863.7273
$ws.Cells.Item(3, 9).Value = 800.1
$ws.Cells.Item(3, 11).Value = 800.1
$ws.Cells.Item(3, 13).Value = -686.1
$ws.Cells.Item(54, 8).Value = 18900
$ws.Cells.Item(54, 9).Value = 18900
$ws.Cells.Item(54, 11).Value = 18900
$ws.Cells.Item(54, 13).Value = -18416
$ws.Cells.Item(86, 8).Value = 3250.3333
$ws.Cells.Item(86, 9).Value = 2253
$ws.Cells.Item(86, 11).Value = 2253
$ws.Cells.Item(86, 13).Value = -1130
$ws.Cells.Item(89, 8).Value = 3250.3333
$ws.Cells.Item(89, 9).Value = 2253
$ws.Cells.Item(89, 11).Value = 11265
$ws.Cells.Item(89, 13).Value = -5649
$ws.Cells.Item(107, 8).Value = 1478
$ws.Cells.Item(107, 9).Value = 1657
$ws.Cells.Item(107, 10).Value = 1299
$ws.Cells.Item(107, 11).Value = 1657
$ws.Cells.Item(107, 12).Value = 1299
$ws.Cells.Item(107, 13).Value = 263
$ws.Cells.Item(107, 14).Value = -5139
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2400.2778
$ws.Cells.Item(31, 9).Value = 1688.7333
$ws.Cells.Item(31, 11).Value = 1688.7333
$ws.Cells.Item(31, 13).Value = -1393.7333
$ws.Cells.Item(34, 8).Value = 2400.2778
$ws.Cells.Item(34, 9).Value = 1688.7333
$ws.Cells.Item(34, 11).Value = 1688.7333
$ws.Cells.Item(34, 13).Value = -1486.7333
$ws.Cells.Item(93, 8).Value = 11981.4
$ws.Cells.Item(93, 9).Value = 11981.4
$ws.Cells.Item(93, 11).Value = 11981.4
$ws.Cells.Item(93, 13).Value = -10109.4
$ws.Cells.Item(105, 8).Value = 1499.3334
$ws.Cells.Item(105, 9).Value = 1249
$ws.Cells.Item(105, 10).Value = 2000
$ws.Cells.Item(105, 11).Value = 1249
$ws.Cells.Item(105, 12).Value = 2000
$ws.Cells.Item(105, 13).Value = 498
$ws.Cells.Item(105, 14).Value = -5494
$ws.Cells.Item(107, 8).Value = 533.8
$ws.Cells.Item(107, 9).Value = 224.66667
$ws.Cells.Item(107, 11).Value = 224.66667
$ws.Cells.Item(107, 13).Value = 1695.33333
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 1443.1351
$ws.Cells.Item(4, 9).Value = 1352.6316
$ws.Cells.Item(4, 10).Value = 1538.6666
$ws.Cells.Item(4, 11).Value = 4057.8948
$ws.Cells.Item(4, 12).Value = 4615.9998
$ws.Cells.Item(4, 13).Value = -3945.8948
$ws.Cells.Item(4, 14).Value = -4839.9998
$ws.Cells.Item(5, 8).Value = 2708.3
$ws.Cells.Item(5, 9).Value = 4631
$ws.Cells.Item(5, 10).Value = 1884.2858
$ws.Cells.Item(5, 11).Value = 13893
$ws.Cells.Item(5, 12).Value = 5652.857400000001
$ws.Cells.Item(5, 13).Value = -13781
$ws.Cells.Item(5, 14).Value = -5876.857400000001
$ws.Cells.Item(37, 8).Value = 99672.625
$ws.Cells.Item(37, 10).Value = 99672.625
$ws.Cells.Item(37, 12).Value = 299017.875
$ws.Cells.Item(37, 14).Value = -299241.875
$ws.Cells.Item(104, 8).Value = 1399.5
$ws.Cells.Item(104, 10).Value = 1399
$ws.Cells.Item(104, 12).Value = 4197
$ws.Cells.Item(104, 14).Value = -9439
$ws.Cells.Item(107, 8).Value = 719.6
$ws.Cells.Item(107, 9).Value = 599.8
$ws.Cells.Item(107, 11).Value = 1799.4
$ws.Cells.Item(107, 13).Value = 120.6000000000001
$ws.Cells.Item(115, 8).Value = 0
$ws.Cells.Item(115, 9).Value = 0
$ws.Cells.Item(115, 11).Value = 0
$ws.Cells.Item(115, 13).ClearContents()
$ws.Cells.Item(132, 8).Value = 3283.8572
$ws.Cells.Item(132, 9).Value = 1831.1666
$ws.Cells.Item(132, 11).Value = 16480.4994
$ws.Cells.Item(132, 13).Value = -13950.4994
$ws.Cells.Item(135, 8).Value = 2708.3
$ws.Cells.Item(135, 9).Value = 4631
$ws.Cells.Item(135, 10).Value = 1884.2858
$ws.Cells.Item(135, 11).Value = 41679
$ws.Cells.Item(135, 12).Value = 16958.5722
$ws.Cells.Item(135, 13).Value = -39144
$ws.Cells.Item(135, 14).Value = -22028.5722
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 78.77778000000001
$ws.Cells.Item(2, 9).Value = 31
$ws.Cells.Item(2, 10).Value = 174.33333
$ws.Cells.Item(2, 11).Value = 31
$ws.Cells.Item(2, 12).Value = 174.33333
$ws.Cells.Item(2, 13).Value = 82
$ws.Cells.Item(2, 14).Value = -400.33333
$ws.Cells.Item(63, 8).Value = 50051.5
$ws.Cells.Item(63, 9).Value = 50103
$ws.Cells.Item(63, 10).Value = 50000
$ws.Cells.Item(63, 11).Value = 50103
$ws.Cells.Item(63, 12).Value = 50000
$ws.Cells.Item(63, 13).Value = -49417
$ws.Cells.Item(63, 14).Value = -51372
$ws.Cells.Item(66, 8).Value = 50051.5
$ws.Cells.Item(66, 9).Value = 50103
$ws.Cells.Item(66, 10).Value = 50000
$ws.Cells.Item(66, 11).Value = 150309
$ws.Cells.Item(66, 12).Value = 150000
$ws.Cells.Item(66, 13).Value = -146877
$ws.Cells.Item(66, 14).Value = -156864
$ws.Cells.Item(80, 8).Value = 1786.2858
$ws.Cells.Item(80, 9).Value = 1599.8
$ws.Cells.Item(80, 11).Value = 1599.8
$ws.Cells.Item(80, 13).Value = -601.8
$ws.Cells.Item(83, 8).Value = 1786.2858
$ws.Cells.Item(83, 9).Value = 1599.8
$ws.Cells.Item(83, 11).Value = 7999
$ws.Cells.Item(83, 13).Value = -3007
$ws.Cells.Item(107, 8).Value = 219.5
$ws.Cells.Item(107, 9).Value = 219.5
$ws.Cells.Item(107, 11).Value = 219.5
$ws.Cells.Item(107, 13).Value = 1700.5
$ws.Cells.Item(113, 8).Value = 1816.8334
$ws.Cells.Item(113, 9).Value = 227.75
$ws.Cells.Item(113, 11).Value = 227.75
$ws.Cells.Item(113, 13).Value = 1942.25
$ws.Cells.Item(126, 8).Value = 5676.3076
$ws.Cells.Item(126, 9).Value = 7665.778
$ws.Cells.Item(126, 10).Value = 1200
$ws.Cells.Item(126, 11).Value = 22997.334
$ws.Cells.Item(126, 12).Value = 3600
$ws.Cells.Item(126, 13).Value = -20527.334
$ws.Cells.Item(126, 14).Value = -8540
$ws.Cells.Item(135, 8).Value = 0
$ws.Cells.Item(135, 10).Value = 0
$ws.Cells.Item(135, 12).Value = 0
$ws.Cells.Item(135, 14).ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 17578.889
$ws.Cells.Item(7, 9).Value = 17578.889
$ws.Cells.Item(7, 11).Value = 17578.889
$ws.Cells.Item(7, 13).Value = -17466.889
$ws.Cells.Item(61, 8).Value = 5992.8
$ws.Cells.Item(61, 10).Value = 6656.3335
$ws.Cells.Item(61, 12).Value = 6656.3335
$ws.Cells.Item(61, 14).Value = -7060.3335
$ws.Cells.Item(95, 8).Value = 28800
$ws.Cells.Item(95, 10).Value = 28800
$ws.Cells.Item(95, 12).Value = 28800
$ws.Cells.Item(95, 14).Value = -34292
$ws.Cells.Item(113, 8).Value = 5992.8
$ws.Cells.Item(113, 10).Value = 6656.3335
$ws.Cells.Item(113, 12).Value = 6656.3335
$ws.Cells.Item(113, 14).Value = -10996.3335
$ws.Cells.Item(126, 8).Value = 17578.889
$ws.Cells.Item(126, 9).Value = 17578.889
$ws.Cells.Item(126, 11).Value = 52736.667
$ws.Cells.Item(126, 13).Value = -50266.667
$ws.Cells.Item(130, 8).Value = 24429
$ws.Cells.Item(130, 10).Value = 24429
$ws.Cells.Item(130, 12).Value = 24429
$ws.Cells.Item(130, 14).Value = -34469
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 1086.5
$ws.Cells.Item(100, 9).Value = 1163.4286
$ws.Cells.Item(100, 10).Value = 548
$ws.Cells.Item(100, 11).Value = 2326.8572
$ws.Cells.Item(100, 12).Value = 1096
$ws.Cells.Item(100, 13).Value = -1785.8572
$ws.Cells.Item(100, 14).Value = -2178
